# Update the run manager sheet for the "MISSING_MOCK_DATA" response-field check
# and move the active selection to F3, as part of switching to yml based configuration.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# F2 previously held "code=MOCK_DATA_NOT_SET" - replace it with the new expected code.
$ws.Range("F2").Value = "code=MISSING_MOCK_DATA"

# Move/update the active selection on the sheet to F3.
$ws.Range("F3").Select()
